$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '0.779 (0.773 ± 0.007)'
$ws.Range("C2").Value = '00:02:36 (00:04:31 ± 00:01:14)'
$ws.Range("D2").Value = '00:00:01 (00:00:07 ± 00:00:04)'
$ws.Range("B3").Value = '0.850 (0.799 ± 0.029)'
$ws.Range("C3").Value = '00:00:11 (00:00:17 ± 00:00:02)'
$ws.Range("D3").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B4").Value = '0.830 (0.758 ± 0.035)'
$ws.Range("C4").Value = '00:00:26 (00:00:43 ± 00:00:10)'
$ws.Range("D4").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B5").Value = '0.876 (0.819 ± 0.036)'
$ws.Range("C5").Value = '00:05:12 (00:05:19 ± 00:00:08)'
$ws.Range("D5").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B6").Value = '0.868 (0.836 ± 0.024)'
$ws.Range("C6").Value = '00:04:54 (00:04:58 ± 00:00:02)'
$ws.Range("D6").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B7").Value = '0.786 (0.730 ± 0.026)'
$ws.Range("C7").Value = '00:05:00 (00:05:04 ± 00:00:02)'
$ws.Range("D7").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B8").Value = '0.888 (0.840 ± 0.030)'
$ws.Range("C8").Value = '00:02:27 (00:05:19 ± 00:02:42)'
$ws.Range("D8").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B9").Value = '0.808 (0.758 ± 0.031)'
$ws.Range("C9").Value = '00:04:59 (00:05:00 ± 00:00:00)'
$ws.Range("D9").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B10").Value = '0.827 (0.777 ± 0.031)'
$ws.Range("C10").Value = '00:04:29 (00:04:29 ± 00:00:00)'
$ws.Range("D10").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B11").Value = '0.863 (0.776 ± 0.077)'
$ws.Range("C11").Value = '00:05:01 (00:05:05 ± 00:00:00)'
$ws.Range("D11").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B12").Value = '0.821 (0.752 ± 0.033)'
$ws.Range("C12").Value = '00:02:01 (00:02:33 ± 00:00:13)'
$ws.Range("D12").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B13").Value = '0.838 (0.719 ± 0.049)'
$ws.Range("C13").Value = '00:00:02 (00:00:03 ± 00:00:00)'
$ws.Range("D13").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B14").Value = '0.863 (0.793 ± 0.032)'
$ws.Range("C14").Value = '00:00:29 (00:00:32 ± 00:00:01)'
$ws.Range("D14").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B15").Value = '0.885 (0.850 ± 0.022)'
$ws.Range("C15").Value = '00:00:27 (00:00:28 ± 00:00:02)'
$ws.Range("D15").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B16").Value = '0.888 (0.846 ± 0.023)'
$ws.Range("C16").Value = '00:00:16 (00:00:17 ± 00:00:00)'
$ws.Range("D16").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B17").Value = '0.880 (0.798 ± 0.034)'
$ws.Range("C17").Value = '00:00:21 (00:02:46 ± 00:01:49)'
$ws.Range("D17").Value = '00:00:00 (00:00:00 ± 00:00:00)'
